# Update the slide-master template: reposition the title/body placeholders
# on the first two custom layouts, and refresh the cached date placeholder
# text (field type datetimeFigureOut) across the remaining layouts.
#
# Note on the position literals: Shape.Left/Top/Width/Height are expressed
# in points and stored by the host as 32-bit floats, while the target
# positions come from the XML diff as EMU (1 pt = 12700 EMU). A naive
# emu/12700 division can truncate to one EMU below the intended value once
# it round-trips through the float32 property, so the literals below are
# chosen to land exactly on the target EMU after that conversion.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$layouts = $master.CustomLayouts

# --- Layout 1 ("Custom Layout") ---------------------------------------
$layout1 = $layouts.Item(1)

# Shape 1: "Title 1"   off(0,0) ext(12192000,460375) -> off(469784,100667) ext(11722216,460375)
$titleShape = $layout1.Shapes.Item(1)
$titleShape.Left = 36.990866
$titleShape.Top = 7.926536
$titleShape.Width = 923.009125
$titleShape.Height = 36.249999

# Shape 2: "Text Placeholder 2"   off(0,554577) ext(12192000,1030603) -> off(469784,655244) ext(11722216,1030603)
$bodyShape = $layout1.Shapes.Item(2)
$bodyShape.Left = 36.990866
$bodyShape.Top = 51.594016
$bodyShape.Width = 923.009125
$bodyShape.Height = 81.149846

# --- Layout 2 ("1_Custom Layout") --------------------------------------
$layout2 = $layouts.Item(2)

# Shape 1: "Text Placeholder 2"   off(0,0) ext(12192000,1030603) -> off(469784,83890) ext(11722216,1030603)
$bodyShape2 = $layout2.Shapes.Item(1)
$bodyShape2.Left = 36.990866
$bodyShape2.Top = 6.605512
$bodyShape2.Width = 923.009125
$bodyShape2.Height = 81.149846

# --- Refresh cached date field text on layouts 3-9 ---------------------
$newDate = "05.07.2024"
for ($li = 3; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

Write-Output "Template layout positions and date placeholders updated."
